$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting old D:K to F:M
$ws.Range("D:E").Insert()

# Copy number/date formatting from column F (the shifted former column D) into new D:E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate new column D and E with the newest two quarters of data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1968000
$ws.Range("E8").Value = 2008000
$ws.Range("D9").Value = 600000
$ws.Range("E9").Value = 655000
$ws.Range("D10").Value = 1368000
$ws.Range("E10").Value = 1353000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 4000
$ws.Range("E14").Value = 112000
$ws.Range("D15").Value = 129000
$ws.Range("E15").Value = 125000
$ws.Range("D17").Value = 1747000
$ws.Range("E17").Value = 1925000
$ws.Range("D18").Value = 221000
$ws.Range("E18").Value = 83000
$ws.Range("D20").Value = 2000
$ws.Range("E20").Value = 6000
$ws.Range("D21").Value = 352000
$ws.Range("E21").Value = 214000
$ws.Range("D22").Value = 22000
$ws.Range("E22").Value = 21000
$ws.Range("D23").Value = 201000
$ws.Range("E23").Value = 68000
$ws.Range("D24").Value = 60000
$ws.Range("E24").Value = 18000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 141000
$ws.Range("E26").Value = 50000
$ws.Range("D27").Value = 141000
$ws.Range("E27").Value = 50000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 28000
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -2000
$ws.Range("E32").Value = -6000
$ws.Range("D33").Value = 169000
$ws.Range("E33").Value = 50000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 169000
$ws.Range("E35").Value = 50000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 474000
$ws.Range("E41").Value = 454000
$ws.Range("D42").Value = 413000
$ws.Range("E42").Value = 483000
$ws.Range("D43").Value = 211000
$ws.Range("E43").Value = 243000
$ws.Range("D44").Value = 78000
$ws.Range("E44").Value = 75000
$ws.Range("D45").Value = 298000
$ws.Range("E45").Value = 252000
$ws.Range("D46").Value = 1474000
$ws.Range("E46").Value = 1507000
$ws.Range("D47").Value = 3000
$ws.Range("E47").Value = 1000
$ws.Range("D48").Value = 8315000
$ws.Range("E48").Value = 8008000
$ws.Range("D49").Value = 96000
$ws.Range("E49").Value = "NA"
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 538000
$ws.Range("E52").Value = 631000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 10426000
$ws.Range("E54").Value = 10147000
$ws.Range("D57").Value = 437000
$ws.Range("E57").Value = 490000
$ws.Range("D58").Value = 309000
$ws.Range("E58").Value = 278000
$ws.Range("D59").Value = 1672000
$ws.Range("E59").Value = 1744000
$ws.Range("D60").Value = 2418000
$ws.Range("E60").Value = 2512000
$ws.Range("D61").Value = 1361000
$ws.Range("E61").Value = 1290000
$ws.Range("D62").Value = 2036000
$ws.Range("E62").Value = 1931000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 5815000
$ws.Range("E66").Value = 5733000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 3679000
$ws.Range("E72").Value = 3511000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 4611000
$ws.Range("E76").Value = 4414000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 169000
$ws.Range("E81").Value = 50000
$ws.Range("D83").Value = 129000
$ws.Range("E83").Value = 125000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 281000
$ws.Range("E89").Value = 206000
$ws.Range("D91").Value = -83000
$ws.Range("E91").Value = -85000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -381000
$ws.Range("E94").Value = -434000
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 120000
$ws.Range("E100").Value = 79000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 20000
$ws.Range("E102").Value = -149000

# Apply individual data corrections in the shifted historical columns (F:M)
$ws.Range("H8").Value = 1758000
$ws.Range("I8").Value = 1818000
$ws.Range("H10").Value = 1261000
$ws.Range("I10").Value = 1337000
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("H17").Value = 1568000
$ws.Range("I17").Value = 1504000
$ws.Range("H18").Value = 190000
$ws.Range("I18").Value = 314000
$ws.Range("H21").Value = 309000
$ws.Range("I21").Value = 431000
$ws.Range("H23").Value = 171000
$ws.Range("I23").Value = 297000
$ws.Range("H24").Value = 67000
$ws.Range("I24").Value = 116000
$ws.Range("H26").Value = 104000
$ws.Range("I26").Value = 181000
$ws.Range("H27").Value = 104000
$ws.Range("I27").Value = 181000
$ws.Range("H29").Value = 551000
$ws.Range("H33").Value = 655000
$ws.Range("I33").Value = 181000
$ws.Range("H35").Value = 655000
$ws.Range("I35").Value = 181000
$ws.Range("H81").Value = 655000
$ws.Range("I81").Value = 181000
$ws.Range("F91").Value = -19000
$ws.Range("G91").Value = -19000
$ws.Range("H91").Value = -38000
$ws.Range("I91").Value = -25000
$ws.Range("J91").Value = -25000
